# Updated symbol list on Fri Jan 13 15:55:02 UTC 2023 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns for the crypto rows
# that moved since the last snapshot. Values are written as text (not
# numbers) so the exact display strings - including trailing zeros and
# the trailing "%" sign - are preserved just like the scraped source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$updates = @(
    @{ Cell = "D2";  Value = "287.43" },
    @{ Cell = "E2";  Value = "2.05%" },
    @{ Cell = "D3";  Value = "29.38" },
    @{ Cell = "E3";  Value = "4.25%" },
    @{ Cell = "D4";  Value = "5.104" },
    @{ Cell = "E4";  Value = "1.47%" },
    @{ Cell = "D5";  Value = "0.06980" },
    @{ Cell = "E5";  Value = "7.69%" },
    @{ Cell = "D6";  Value = "7.418" },
    @{ Cell = "E6";  Value = "2.39%" },
    @{ Cell = "D7";  Value = "3.575" },
    @{ Cell = "E7";  Value = "5.93%" },
    @{ Cell = "E8";  Value = "-0.83%" },
    @{ Cell = "D9";  Value = "0.9026" },
    @{ Cell = "E9";  Value = "-2.92%" },
    @{ Cell = "D10"; Value = "0.1595" },
    @{ Cell = "E10"; Value = "3.04%" },
    @{ Cell = "D11"; Value = "0.07123" },
    @{ Cell = "E11"; Value = "16.25%" },
    @{ Cell = "D12"; Value = "0.07675" },
    @{ Cell = "E12"; Value = "2.08%" },
    @{ Cell = "D13"; Value = "0.02911" },
    @{ Cell = "E13"; Value = "0.23%" },
    @{ Cell = "D14"; Value = "0.08988" },
    @{ Cell = "E14"; Value = "0.21%" },
    @{ Cell = "D15"; Value = "0.001595" },
    @{ Cell = "E15"; Value = "0.77%" },
    @{ Cell = "D16"; Value = "0.0006482" },
    @{ Cell = "E16"; Value = "1.22%" },
    @{ Cell = "D17"; Value = "0.006404" },
    @{ Cell = "E17"; Value = "5.96%" },
    @{ Cell = "E18"; Value = "0.43%" },
    @{ Cell = "D19"; Value = "2.230" },
    @{ Cell = "E19"; Value = "-0.14%" },
    @{ Cell = "D20"; Value = "0.3234" },
    @{ Cell = "E20"; Value = "1.31%" },
    @{ Cell = "E21"; Value = "1.50%" },
    @{ Cell = "D22"; Value = "3.996" },
    @{ Cell = "E22"; Value = "-1.71%" },
    @{ Cell = "D23"; Value = "0.1554" },
    @{ Cell = "E23"; Value = "0.61%" },
    @{ Cell = "D24"; Value = "0.04507" },
    @{ Cell = "E24"; Value = "1.49%" },
    @{ Cell = "D25"; Value = "0.001205" },
    @{ Cell = "E25"; Value = "1.92%" },
    @{ Cell = "D26"; Value = "0.004384" },
    @{ Cell = "E26"; Value = "-0.21%" },
    @{ Cell = "D27"; Value = "0.0001165" },
    @{ Cell = "E27"; Value = "-6.72%" },
    @{ Cell = "D28"; Value = "0.0001612" },
    @{ Cell = "E28"; Value = "-0.33%" },
    @{ Cell = "D40"; Value = "0.04286" },
    @{ Cell = "E40"; Value = "3.25%" },
    @{ Cell = "D41"; Value = "0.006815" },
    @{ Cell = "E41"; Value = "2.90%" },
    @{ Cell = "D42"; Value = "0.1248" },
    @{ Cell = "E42"; Value = "2.24%" },
    @{ Cell = "D43"; Value = "0.002211" },
    @{ Cell = "E43"; Value = "9.53%" },
    @{ Cell = "D44"; Value = "0.01152" },
    @{ Cell = "E44"; Value = "-4.27%" },
    @{ Cell = "D45"; Value = "0.00005744" },
    @{ Cell = "E45"; Value = "2.36%" },
    @{ Cell = "E46"; Value = "-1.85%" },
    @{ Cell = "D47"; Value = "0.01302" },
    @{ Cell = "E47"; Value = "0.21%" }
)

foreach ($u in $updates) {
    Set-TextValue $u.Cell $u.Value
}
